# agrega campo Activo a tabla Usuarios
# -------------------------------------------------------------
# This script reproduces the authored change: the "Usuario" table
# gets its header row re-cased to UPPERCASE and a new "ACTIVO" flag
# column (F) is appended, populated with 1 for every active user and
# 0 for the last (inactive) row. It also re-applies the sheet tab
# colors and the active-sheet/selection bookkeeping that Excel wrote
# when the workbook was re-saved.
# -------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- Usuario sheet: header text + new Activo column ----------
$wsUsuario = $wb.Worksheets.Item("Usuario")

$wsUsuario.Range("A1").Value = "IDUSUARIO"
$wsUsuario.Range("B1").Value = "USUARIO"
$wsUsuario.Range("C1").Value = "TIPOUSUARIO"
$wsUsuario.Range("D1").Value = "CLAVE"
$wsUsuario.Range("E1").Value = "MAIL"
$wsUsuario.Range("F1").Value = "ACTIVO"

# Rows 2-10 are active users, row 11 (the last Prevencionista) is not
$wsUsuario.Range("F2:F10").Value = 1
$wsUsuario.Range("F11").Value = 0

# ---- Sheet tab colors --------------------------------------------------
# theme Accent4 (theme="7") -> RGB FFC000
$wb.Worksheets.Item("Usuario").Tab.Color = 49407
$wb.Worksheets.Item("Admin").Tab.Color = 49407
$wb.Worksheets.Item("Profesional").Tab.Color = 49407
$wb.Worksheets.Item("Cliente").Tab.Color = 49407

# explicit rgb colors
$wb.Worksheets.Item("ReporteAccidente").Tab.Color = 10498160   # FF7030A0
$wb.Worksheets.Item("Actividad").Tab.Color = 15773696          # FF00B0F0
$wb.Worksheets.Item("Factura").Tab.Color = 5287936             # FF00B050
$wb.Worksheets.Item("Item").Tab.Color = 5287936                # FF00B050

# theme Accent2 (theme="5"), tint -0.249977111117893 ("Darker 25%") -> RGB C55A11
$wb.Worksheets.Item("SolicitudAsesoria").Tab.Color = 1137349

# ---- Selections / active sheet ----------------------------------------
# Usuario: selection moves off the old A4:A7 block
$wsUsuario.Activate()
$wsUsuario.Range("B30").Select()

# Item: selection moves, and it gives up being the active tab
$wsItem = $wb.Worksheets.Item("Item")
$wsItem.Activate()
$wsItem.Range("G17").Select()

# Factura becomes the active tab (tabSelected) in its place
$wb.Worksheets.Item("Factura").Activate()
